# Update "metricas_retencao_anual" data: bump num_customers (and derived
# cohort_size/retention_rate) for a handful of rows to reflect updated BIBI data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: cohort 2021, period 4 -> num_customers 56 -> 57
$ws.Range("C27").Value = 57
$ws.Range("E27").Value = 57 / 2252

# Row 31: cohort 2022, period 3 -> num_customers 56 -> 57
$ws.Range("C31").Value = 57
$ws.Range("E31").Value = 57 / 2312

# Row 34: cohort 2023, period 2 -> num_customers 89 -> 90
$ws.Range("C34").Value = 90
$ws.Range("E34").Value = 90 / 2256

# Row 36: cohort 2024, period 1 -> num_customers 145 -> 146
$ws.Range("C36").Value = 146
$ws.Range("E36").Value = 146 / 1930

# Row 37: cohort 2025, period 0 -> num_customers & cohort_size 955 -> 967
$ws.Range("C37").Value = 967
$ws.Range("D37").Value = 967
